# Apply the "cryptos list" refresh described by the commit diff.
# D-column prices are plain text (e.g. "29.353.53", "1.000") that must
# NOT be auto-coerced into numbers by Excel, so NumberFormat is forced
# to Text ("@") before writing each D-column value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.353.53'

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.842.66'
$ws.Range("E3").Value = '  -0.81%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.10'
$ws.Range("E5").Value = '  -0.68%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6294'
$ws.Range("E6").Value = '  -0.64%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9996'
$ws.Range("E7").Value = '  -0.32%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07399'
$ws.Range("E8").Value = '  -2.44%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2900'
$ws.Range("E9").Value = '  -0.99%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.85'
$ws.Range("E10").Value = '  +0.91%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07738'
$ws.Range("E11").Value = '  -0.36%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.840.22'
$ws.Range("E12").Value = '  -0.93%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.981'
$ws.Range("E13").Value = '  -1.28%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6797'
$ws.Range("E14").Value = '  -0.92%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001022'
$ws.Range("E15").Value = '  -2.41%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '82.01'

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.259'
$ws.Range("E17").Value = '  +1.62%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '29.346.92'
$ws.Range("E18").Value = '  -0.50%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '229.25'
$ws.Range("E19").Value = '  -0.50%  '

$ws.Range("E20").Value = '  -0.64%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9999'
$ws.Range("E21").Value = '  -0.29%  '

$ws.Range("E22").Value = '  -1.09%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.000'
$ws.Range("E23").Value = '  -0.34%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '158.24'
$ws.Range("E24").Value = '  -0.86%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.479'

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1355'
$ws.Range("E26").Value = '  -3.28%  '

$ws.Range("E27").Value = '  -1.77%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.06526'
$ws.Range("E28").Value = '  +14.31%  '

$ws.Range("E29").Value = '  +2.04%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.488'
$ws.Range("E30").Value = '  +0.38%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.073'
$ws.Range("E31").Value = '  -2.07%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.065'
$ws.Range("E32").Value = '  -0.07%  '

$ws.Range("E33").Value = '  +0.11%  '

$ws.Range("E34").Value = '  -1.79%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6941'
$ws.Range("E35").Value = '  -0.21%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.573'
$ws.Range("E36").Value = '  -0.92%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01859'
$ws.Range("E37").Value = '  +1.31%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.818'
$ws.Range("E38").Value = '  +1.18%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.242.19'
$ws.Range("E39").Value = '  -1.18%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.807'
$ws.Range("E40").Value = '  +4.38%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9338'
$ws.Range("E41").Value = '  +2.78%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9992'
$ws.Range("E42").Value = '  -0.36%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.988.00'
$ws.Range("E43").Value = '  -1.51%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '100.87'
$ws.Range("E44").Value = '  -0.88%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '65.59'
$ws.Range("E45").Value = '  -1.04%  '

$ws.Range("E46").Value = '  +4.07%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.060'
$ws.Range("E47").Value = '  -1.53%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.714'
$ws.Range("E48").Value = '  +2.02%  '

$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.1153'
$ws.Range("E49").Value = '  -1.61%  '

$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.995'
$ws.Range("E50").Value = '  -0.55%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3905'
$ws.Range("E51").Value = '  -1.96%  '
